$wb = $excel.ActiveWorkbook

# --- MigrationPage: just move the selection cursor (no data change) ---
$wsMigration = $wb.Worksheets.Item("MigrationPage")
$wsMigration.Range("B30").Select() | Out-Null

# --- PlansAndAddonsPage: add the two new xpath rows (AddOnsSelect / Subcribe) ---
$wsPlans = $wb.Worksheets.Item("PlansAndAddonsPage")

$rowAddOnsSelectKey = $wsPlans.Range("A8")
$rowAddOnsSelectVal = $wsPlans.Range("B8")
$rowAddOnsSelectKey.Value = "AddOnsSelect"
$rowAddOnsSelectVal.Value = "//android.widget.LinearLayout[1]/android.view.ViewGroup"
$rowAddOnsSelectVal.Style = "Normal 7"
$rowAddOnsSelectVal.WrapText = $true

$rowSubscribeKey = $wsPlans.Range("A9")
$rowSubscribeVal = $wsPlans.Range("B9")
$rowSubscribeKey.Value = "Subcribe"
$rowSubscribeVal.Value = "//android.widget.TextView[@text='SUBSCRIBE']"
$rowSubscribeVal.Style = "Normal 7"
$rowSubscribeVal.WrapText = $true

# restore this sheet (PlansAndAddonsPage) as the active tab/selection, matching the diff
$wsPlans.Range("B17").Select() | Out-Null
